# Saldo.xlsx update ("Add files via upload")
#
# The refreshed export:
#  - drops the "004450724 / Assako / 301650.6" row
#  - replaces the "004479463 / Henrique / 29208.11" row with
#    "004453157 / Jose / 25862.99", followed by a brand-new row
#    "002064834 / Rafaela / 2984.7"
#  - replaces the "004556853 / Marcel / 2502.95" row with
#    "008243633 / Daniela / 2400"
#  - drops the "005685089 / Carneiro / 155.97" row
#
# Row numbers below refer to the *original* sheet layout (header is row 1).
# Operations are applied from the bottom of the sheet upward so that each
# target row number is still valid when it is used (earlier deletes/inserts
# further down the sheet would otherwise shift the row indices of the
# entries still waiting to be processed above them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "005685089 / Carneiro / 155.97" row entirely.
$ws.Rows(119).Delete()

# Turn "004556853 / Marcel / 2502.95" into "008243633 / Daniela / 2400".
# Format column A as Text first so the leading zero in the account number
# is preserved instead of Excel interpreting it as a number.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "008243633"
$ws.Range("B7").Value = "Daniela"
$ws.Range("C7").Value = 2400

# Turn "004479463 / Henrique / 29208.11" into "004453157 / Jose / 25862.99".
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "004453157"
$ws.Range("B5").Value = "Jose"
$ws.Range("C5").Value = 25862.99

# Insert a brand-new row right after it for "002064834 / Rafaela / 2984.7".
$ws.Rows(6).Insert()
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "002064834"
$ws.Range("B6").Value = "Rafaela"
$ws.Range("C6").Value = 2984.7

# Remove the "004450724 / Assako / 301650.6" row entirely.
$ws.Rows(2).Delete()
